$d = $word.ActiveDocument

function Find-ParagraphIndex($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# -------------------------------------------------------------------
# Hunk 1: replace the lone-space paragraph (end of the W006 cell) with
# an empty paragraph followed by two new paragraphs ("5 nov 2020" and
# the long status update ending in a 3-run split around "nog").
# -------------------------------------------------------------------

$idx = Find-ParagraphIndex("Doel is nog om de optieprijs")
$spaceIdx = $idx + 1
$spacePara = $d.Paragraphs.Item($spaceIdx)
$spaceRange = $spacePara.Range

# Drop the run text (" ") but keep the paragraph mark & pPr so the
# paragraph becomes empty, matching the diff's removed <w:r>.
$clearRange = $d.Range($spaceRange.Start, $spaceRange.End - 1)
$clearRange.Delete()

# Insert the two following paragraphs (formatting is inherited from
# the now-empty paragraph, which already carries the nl-NL lang rPr).
$spacePara2 = $d.Paragraphs.Item($spaceIdx)
$insertRange = $spacePara2.Range
$insertRange.InsertAfter("`r5 nov 2020`r")

$longIdx = $spaceIdx + 2
$longPara = $d.Paragraphs.Item($longIdx)
$longRange = $longPara.Range
$longRange.InsertAfter("Optieprijs kan nog niet worden aangepast, optietransactie wordt wel netjes ingevoerd aan de hand van optiepremie en contractgrootte kan 100 zijn of 10. Totalen kloppen (eindelijk) ook weer. Inlezen en naar bestand schrijven van posities en transacties netjes naar de objecten verplaatst, de intelligentie rond berekenen van totalen zwerft nog een beetje rond en is nog aanwezig in het frontend, moet nog verplaatst worden naar de positie- en transactieobjecten.")

# Split "moet nog verplaatst" -> "...moet " | "nog" | " verplaatst..."
# so the run boundaries match the diff (Bold toggle forces a run
# break without leaving stray formatting behind).
$longPara2 = $d.Paragraphs.Item($longIdx)
$fullRange = $longPara2.Range
$narrow = $d.Range($fullRange.Start, $fullRange.End)
[void]$narrow.Find.Execute("moet nog verplaatst")
$nogRange = $d.Range($narrow.Start, $narrow.End)
[void]$nogRange.Find.Execute("nog")
$nogRange.Bold = 1
$nogRange.Bold = 0
